$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..AY hold data for rows 2..13 (row 1 is the header).
# This edit relocates each data row's full content to a different row
# (a permutation of rows 2-13), as derived from the source diff:
#   old row -> new row
#   2  -> 4
#   3  -> 5
#   4  -> 7
#   5  -> 13
#   6  -> 2
#   7  -> 3
#   8  -> 6
#   9  -> 8
#   10 -> 9
#   11 -> 10
#   12 -> 11
#   13 -> 12
$mapping = @{2=4; 3=5; 4=7; 5=13; 6=2; 7=3; 8=6; 9=8; 10=9; 11=10; 12=11; 13=12}

$firstRow = 2
$lastRow = 13
$lastCol = 51   # column AY

# --- Pass 1: snapshot every source cell's value before anything is overwritten ---
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

# --- Pass 2: clear the old data area so stale values can't linger ---
$ws.Range("A2:AY13").ClearContents()

# --- Pass 3: write every cell to its new row according to the mapping ---
foreach ($oldRow in $mapping.Keys) {
    $newRow = $mapping[$oldRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $v = $snapshot["$oldRow-$c"]
        if ($v -eq $null) {
            continue
        }
        # Empty-string placeholders (blank inlineStr cells) carry no data;
        # skip them so the cell is simply left empty. NOTE: don't use
        # "$v -eq ''' for this check — PowerShell coerces $false -eq ''
        # to $true, which would wrongly drop boolean FALSE values.
        $typeName = $v.GetType().Name
        if ($typeName -eq "String" -and $v -eq "") {
            continue
        }
        $target = $ws.Cells.Item($newRow, $c)
        if ($typeName -eq "String") {
            # Force text so date/time-looking strings (e.g. "2016-07-05",
            # "00:00") are not auto-coerced into date/time serials.
            $target.Value = "'" + $v
        } else {
            # Numbers and booleans round-trip safely as-is.
            $target.Value = $v
        }
    }
}
